$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Future Heads" helper table in columns N:O (rows 2-20) is no longer
# needed now that the logic grabs food immediately when it is close by,
# so clear that whole helper area. Clearing (rather than deleting rows/
# columns) matches the original layout: some rows only had content in
# N:O and disappear once empty, while rows 16-19 keep their B:C content.
$ws.Range("N2:O20").ClearContents()

# Reflect the new selection left behind after the edit.
$ws.Range("M2").Select()
